# Commit: "changing document, table attributes to lowerCamelCase"
#
# The ObjTables header rows embedded as text in the first row(s) of each
# sheet use a set of pseudo-XML-ish attributes (Type=, Id=, Description=,
# Name=, Date=, ObjTablesVersion=) that need to be renamed to
# lowerCamelCase (type=, id=, description=, name=, date=,
# objTablesVersion=) while keeping the rest of the row untouched.

$wb = $excel.ActiveWorkbook

# "!!_Table of contents" sheet: row 1 (A1) and row 2 (A2) headers.
$wsToc = $wb.Worksheets.Item("!!_Table of contents")
$wsToc.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8'"
$wsToc.Range("A2").Value = "!!ObjTables type='TableOfContents' description='Table of contents' date='2019-09-22 23:16:33' objTablesVersion='0.0.8'"

# "!!_Schema" sheet: row 1 (A1) header.
$wsSchema = $wb.Worksheets.Item("!!_Schema")
$wsSchema.Range("A1").Value = "!!ObjTables type='Schema' description='Table/model and column/attribute definitions' date='2019-09-22 23:16:33' objTablesVersion='0.0.8'"

# "!!Compound" sheet: row 1 (A1) header.
$wsCompound = $wb.Worksheets.Item("!!Compound")
$wsCompound.Range("A1").Value = "!!ObjTables type='Data' id='Compound' description='Compound' name='Compound' date='2019-09-22 23:16:33' objTablesVersion='0.0.8'"

# "!!Model" sheet: row 1 (A1) header.
$wsModel = $wb.Worksheets.Item("!!Model")
$wsModel.Range("A1").Value = "!!ObjTables type='Data' id='Model' description='Model' name='Model' date='2019-09-22 23:16:33' objTablesVersion='0.0.8'"

# "!!Reaction" sheet: row 1 (A1) header.
$wsReaction = $wb.Worksheets.Item("!!Reaction")
$wsReaction.Range("A1").Value = "!!ObjTables type='Data' id='Reaction' description='Reaction' name='Reaction' date='2019-09-22 23:16:33' objTablesVersion='0.0.8'"
